$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Row 2: OrchestratorQueueName value changes from "ProcessABCQueue" to "WI5_Items"
$ws.Range("B2").Value = "WI5_Items"

# Row 3: OrchestratorQueueFolder value changes from empty to "Shared"
$ws.Range("B3").Value = "Shared"

# New settings rows added to the Settings sheet
$ws.Range("A7").Value = "System1Url"
$ws.Range("B7").Value = "https://acme-test.uipath.com/"

$ws.Range("A9").Value = "System1WorkItems"
$ws.Range("B9").Value = "https://acme-test.uipath.com/work-items"

$ws.Range("A11").Value = "SHA1OnlineUrl"
$ws.Range("B11").Value = "http://www.sha1-online.com/"

$ws.Range("A13").Value = "System1Credential"
$ws.Range("B13").Value = "System1Credential"

$ws.Range("A15").Value = "ProcessName"
$ws.Range("B15").Value = "chrome"

$ws.Range("A17").Value = "System1ResetDataUrl"
$ws.Range("B17").Value = "https://acme-test.uipath.com/reset-test-data"

# Reflect the cursor position left in the sheet after editing
$ws.Range("B6").Select() | Out-Null

$wb.Save() | Out-Null
